# fix print treeview issue
$wb = $excel.ActiveWorkbook

# Sheet "식당판매" - C7 was blank, set to 0
$ws1 = $wb.Worksheets.Item("식당판매")
$ws1.Range("C7").Value = 0

# Sheet "매점판매" - C2 was blank, set to 0
$ws2 = $wb.Worksheets.Item("매점판매")
$ws2.Range("C2").Value = 0

# Sheet "기타" - update several totals
$ws5 = $wb.Worksheets.Item("기타")
$ws5.Range("C7").Value = 10
$ws5.Range("C8").Value = 31
$ws5.Range("C10").Value = 27
$ws5.Range("C11").Value = 37
